$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 28 first (higher row number) so row indices for the earlier
# deletion stay valid, then remove row 26 ("RM 232").
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()
